$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9642857142857143
$ws.Range("C2").Value = 0.7694805194805194

$ws.Range("B3").Value = 0.9675324675324676
$ws.Range("C3").Value = 0.75

$ws.Range("B4").Value = 0.9707792207792207
$ws.Range("C4").Value = 0.75

$ws.Range("B5").Value = 0.974025974025974
$ws.Range("C5").Value = 0.762987012987013

$ws.Range("B6").Value = 0.9675324675324676
$ws.Range("C6").Value = 0.7305194805194806
